# Apply "additions to metadata table" edit:
# Adds a new row 12 to Sheet1 with a third RNAseq dataset entry:
#   B12 = "Ran dataset 3"
#   C12 = date 12/6/2024 (serial 45632), formatted as a short date
#   D12 = "Index of /emma.strand/rnaseq/Cvir_Prkns_rnaseq_dataset3" with a
#         hyperlink to https://gannet.fish.washington.edu/emma.strand/rnaseq/Cvir_Prkns_rnaseq_dataset3/
# Also updates the active selection in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row ---
$ws.Range("B12").Value = "Ran dataset 3"

$ws.Range("C12").Value = 45632
$ws.Range("C12").NumberFormat = "mm-dd-yy"

$ws.Range("D12").Value = "Index of /emma.strand/rnaseq/Cvir_Prkns_rnaseq_dataset3"

# --- Hyperlink on D12 ---
$link = $ws.Hyperlinks.Add($ws.Range("D12"), "https://gannet.fish.washington.edu/emma.strand/rnaseq/Cvir_Prkns_rnaseq_dataset3/")
$link.TextToDisplay = "https://gannet.fish.washington.edu/emma.strand/rnaseq/Cvir_Prkns_rnaseq_dataset3/"
# restore the cell's visible text (TextToDisplay assignment above only affects
# the hyperlink's stored display text, but set Value again defensively)
$ws.Range("D12").Value = "Index of /emma.strand/rnaseq/Cvir_Prkns_rnaseq_dataset3"

# --- Update selection shown in the sheet view ---
$ws.Range("C22").Select() | Out-Null
